# Add "Homework 3" and "Quiz 1" columns (quiz 1 results) to the grade book.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column headers
$ws.Range("F1").Value = "Homework 3"
$ws.Range("G1").Value = "Quiz 1"

# Row 2 - Gary He
$ws.Range("F2").Formula = "=29/30"
$ws.Range("G2").Formula = "=55/50"

# Row 3 - Jayden Warlum
$ws.Range("E3").Formula = "=24/25"
$ws.Range("F3").Formula = "=29/30"
$ws.Range("G3").Formula = "=44/50"

# Row 4 - Aman Kumpawat
$ws.Range("E4").Formula = "=27/25"
$ws.Range("F4").Formula = "=30/30"
$ws.Range("G4").Formula = "=55/50"

# Row 5 - Nahom Anteneh
$ws.Range("E5").Formula = "=25/25"
$ws.Range("F5").Formula = "=29/30"
$ws.Range("G5").Formula = "=55/50"

# Row 6 - Oswen Martinez
$ws.Range("E6").Value = 0
$ws.Range("F6").Value = 0
$ws.Range("G6").Formula = "=47/50"

# Row 7 - Kai Stephens
$ws.Range("E7").Formula = "=15/25"
$ws.Range("F7").Value = 0
$ws.Range("G7").Formula = "=41/50"

# Row 8 - James Saw
$ws.Range("F8").Formula = "=29/30"
$ws.Range("G8").Formula = "=50/50"

# Row 9 - Edward Whitesel
$ws.Range("F9").Formula = "=28/30"
$ws.Range("G9").Formula = "=44/50"

# Row 10 - Loren Grey
$ws.Range("E10").Formula = "=25/25"
$ws.Range("F10").Formula = "=29/30"
$ws.Range("G10").Formula = "=55/50"

# Row 11 - Almas Waseem
$ws.Range("E11").Formula = "=0"
$ws.Range("F11").Value = 0
$ws.Range("G11").Formula = "=0"

# Row 12 - Ty Carlson
$ws.Range("E12").Value = 0
$ws.Range("F12").Value = 0
$ws.Range("G12").Formula = "=0"

# Row 13 - Nailyn Lopez
$ws.Range("F13").Formula = "=29/30"
$ws.Range("G13").Formula = "=47/50"

# Row 14 - Roy Kalu
$ws.Range("E14").Formula = "=18/25"
$ws.Range("F14").Formula = "=29/30"
$ws.Range("G14").Formula = "=38/50"

# Row 15 - Dylan Zeledon
$ws.Range("F15").Formula = "=5/30"
$ws.Range("G15").Formula = "=0"

# Column widths for the two new columns (F, G) - values chosen so the
# persisted OOXML <col width=.../> lands as close as possible to the
# target widths (10.99 / 11.5) given the engine's character-width
# quantization.
$ws.Columns.Item(6).ColumnWidth = 10.084
$ws.Columns.Item(7).ColumnWidth = 10.666666666666666

# Move / record the selection where the author left off editing.
$ws.Range("G11").Select()
